$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns before the old "Message" column (old H, now becomes M)
$ws.Range("H1:L1").EntireColumn.Insert()

# New columns G:L share the same (narrow) width
$ws.Range("G1:L1").ColumnWidth = 7.73333333333333

# New header cells for the inserted columns (written in the same order the
# shared-string table picked them up in: I, J, K, L, then H)
$ws.Range("I1").Value = "AttachDocumentName"
$ws.Range("J1").Value = "AttachSupportDocuments"
$ws.Range("K1").Value = "AttachSupportDocumentName"
$ws.Range("L1").Value = "ReviewDocument"
$ws.Range("H1").Value = "AttachDocuments"

# Row 2 only had a single coloured placeholder cell (G2); the newly inserted
# columns inherited that placeholder style, so clear all of them away.
$ws.Range("H2:L2").Clear()

# Rows 3 and 4 keep a second placeholder cell in the last inserted column (L)
# - matching the style already copied from G - so only clear H:K there.
$ws.Range("H3:K3").Clear()
$ws.Range("H4:K4").Clear()
